$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Purchase 22-23")

# New data row (row 22) — mirror the formatting of row 2, which uses the
# same style indices (9,3,9,9,9,25) as the target row.
$ws.Range("A2:F2").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A22").Value = 3
$ws.Range("B22").Value = 45102
$ws.Range("C22").Value = "SLH/1416"
$ws.Range("D22").Value = "Shree Laxmi Lighting Hub"
$ws.Range("E22").Value = 1033
$ws.Range("F22").Formula = "=E22"

# Update the view: clear the old scroll position and move the active
# selection to A24 (just below the newly added data).
$ws.Range("A24").Select()
